$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 3
$ws.Range("I3").Value = 2.32
$ws.Range("J3").Value = 3.45
$ws.Range("L3").Value = 2.85
$ws.Range("N3").Value = 7.9
$ws.Range("O3").Value = 1.32
$ws.Range("P3").Value = 2.85
$ws.Range("R3").Value = 1.7
$ws.Range("U3").Value = 1.7
$ws.Range("V3").Value = 1.91
$ws.Range("W3").Value = 8.75
$ws.Range("X3").Value = 15.5
$ws.Range("Y3").Value = 10.5
$ws.Range("Z3").Value = 37
$ws.Range("AB3").Value = 35
$ws.Range("AC3").Value = 8.75
$ws.Range("AE3").Value = 13.5
$ws.Range("AF3").Value = 65
$ws.Range("AG3").Value = 7.7
$ws.Range("AI3").Value = 9
$ws.Range("AK3").Value = 19.5
$ws.Range("AL3").Value = 29
$ws.Range("AM3").Value = 500
$ws.Range("AN3").Value = 4.9
$ws.Range("AO3").Value = 16
$ws.Range("AP3").Value = 22
$ws.Range("AS3").Value = 2.57
$ws.Range("AT3").Value = 6.6
$ws.Range("AU3").Value = 55
$ws.Range("AX3").Value = 18.5
$ws.Range("AZ3").Value = 75
